# The sheet originally had columns:
#   A=Nombre, B=Correo electrónico, C=ID, D=Localización, E=Tipo
# with one data row (Pablo Pinto, pablo@example.com, 59687412O, 36S78W, <type value gone>)
#
# The edit:
#   - removes the "Localización" column entirely
#   - replaces the old single "Tipo" value (36S78W) with two new numeric
#     columns "Latitud" / "Longitud", inserted before "Tipo"
#   - "Tipo" header is kept (now the last header, column F) but no data
#     value is filled in for it on the data row
#
# Resulting layout:
#   A=Nombre, B=Correo electrónico, C=ID, D=Latitud, E=Longitud, F=Tipo

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Correo electrónico"
$ws.Range("C1").Value = "ID"
$ws.Range("D1").Value = "Latitud"
$ws.Range("E1").Value = "Longitud"
$ws.Range("F1").Value = "Tipo"

# Data row
$ws.Range("A2").Value = "Pablo Pinto"
$ws.Range("B2").Value = "pablo@example.com"
$ws.Range("C2").Value = "59687412O"
$ws.Range("D2").Value = 156.26
$ws.Range("E2").Value = -10.265000000000001

$ws.Range("E2").Select()
